# [Kadastro App] Yeni kayit eklendi: 2905
# Adds a new record row to both the "Kayitlar" master sheet and the
# "Erdemli" birim sheet, mirroring the existing rows for that birim.

$wb = $excel.ActiveWorkbook

$kayitNo   = "2905"
$tarih     = "2025-09-08"
$birim     = "Erdemli"
$parsel    = "1"
$is        = "3B"
$personel  = "EMİNE ALANLI KIRCILI (K.Mühendisi), SEVİL SARAÇER (Tekniker)"

function Add-Kayit($ws, $row) {
    # Leading apostrophes force text storage (matches the existing rows,
    # which are all plain text cells rather than numbers/dates) without
    # touching the cell's number format/style.
    $ws.Range("A$row").Value = "'" + $kayitNo
    $ws.Range("B$row").Value = "'" + $tarih
    $ws.Range("C$row").Value = $birim
    $ws.Range("D$row").Value = "'" + $parsel
    $ws.Range("E$row").Value = $is
    $ws.Range("F$row").Value = $personel
}

# "Kayitlar" sheet: existing data occupies rows 1-11, new row goes at 12.
$wsKayitlar = $wb.Worksheets.Item("Kayitlar")
Add-Kayit $wsKayitlar 12

# "Erdemli" sheet: existing data occupies rows 1-10, new row goes at 11.
$wsErdemli = $wb.Worksheets.Item("Erdemli")
Add-Kayit $wsErdemli 11
